$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: Swap the data (columns B..AB) between row 73 and row 74 ---
# Column A (the running index) stays as-is for both rows.
$cols = @("B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB")

foreach ($col in $cols) {
    $addr73 = $col + "73"
    $addr74 = $col + "74"
    $v73 = $ws.Range($addr73).Value2
    $v74 = $ws.Range($addr74).Value2
    $ws.Range($addr73).Value = $v74
    $ws.Range($addr74).Value = $v73
}

# --- Step 2: Append a new match record as row 171 ---
$ws.Range("A171").Value = 169
$ws.Range("B171").Value = 8233998
$ws.Range("C171").Value = "Australia ALeague"
$ws.Range("D171").Value = 45437.28125
$ws.Range("E171").Value = "Central Coast Mariners"
$ws.Range("F171").Value = "Melbourne Victory"
$ws.Range("G171").Value = 3
$ws.Range("H171").Value = 1
$ws.Range("I171").Value = "H"
$ws.Range("J171").Value = 2.3
$ws.Range("K171").Value = 3.4
$ws.Range("L171").Value = 3
$ws.Range("M171").Value = 2.35
$ws.Range("N171").Value = 3.1
$ws.Range("O171").Value = 3.25
$ws.Range("P171").Value = -0.25
$ws.Range("Q171").Value = 2.025
$ws.Range("R171").Value = 1.825
$ws.Range("S171").Value = 2.25
$ws.Range("T171").Value = 2.05
$ws.Range("U171").Value = 1.8
$ws.Range("V171").Value = 1.35
$ws.Range("W171").Value = -1
$ws.Range("X171").Value = -1
$ws.Range("Y171").Value = 1.025
$ws.Range("Z171").Value = -1
$ws.Range("AA171").Value = 1.05
$ws.Range("AB171").Value = -1

# Copy the formatting of the previous last row (170) onto the new row 171
# so that the id column (A) and the date column (D) keep their original
# styles (bold/centered border style and date number format respectively).
$ws.Range("A170").Copy() | Out-Null
$ws.Range("A171").PasteSpecial(-4122) | Out-Null

$ws.Range("D170").Copy() | Out-Null
$ws.Range("D171").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0
